# "results from 10000 iters"
# The HCR (harvest control rule) number used to build scenarios 38-43 moves
# from 12 to 13. Column B on those rows is the shared formula
#   =CONCATENATE("ASS",C,"_HCR",D,"_REC",E,"_INN",F,"_OER",G)
# so updating column D recalculates the B column text automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Folha1")

foreach ($r in 38..43) {
    $ws.Cells.Item($r, 4).Value = 13
}

# Make sure every cached formula result (column B) is refreshed.
$excel.CalculateFullRebuild()

# Update the view: scroll the frozen pane down and move the active selection.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 23
$ws.Range("H46").Select()

$wb.Save()
